$d = $word.ActiveDocument

function Get-ParagraphByText($text) {
    foreach ($p in $d.Paragraphs) {
        if ($p.Range.Text.TrimEnd("`r`a") -eq $text) {
            return $p
        }
    }
    return $null
}

function Set-ParagraphText($oldText, $newText) {
    $p = Get-ParagraphByText $oldText
    if ($p -ne $null) {
        $p.Range.Text = $newText
    }
}

function Remove-ParagraphByText($text) {
    $p = Get-ParagraphByText $text
    if ($p -ne $null) {
        $p.Range.Delete()
    }
}

# 1) Rewrite the five existing "Functional Requirements" user-story bullets
# (direct Range.Text assignment keeps straight apostrophes intact and
# preserves the paragraph's ListBullet style/formatting)
Set-ParagraphText "As a potential traveler, I want to search for flights by origin, destination, and date so that I can find available options." "As a new user, I want to create an account so that I can shop and save my information."
Set-ParagraphText "As a traveler, I want to filter flight search results by price, airline, and number of stops so that I can find the best-suited option." "As a registered user, I want to browse the product catalog so that I can find items I'm interested in."
Set-ParagraphText "As a traveler, I want to select a flight and view detailed information about it, including price, time, and stops, so that I can make an informed decision." "As a user, I want to add items to my shopping cart so that I can purchase multiple products in one transaction."
Set-ParagraphText "As a traveler, I want to securely input my personal and payment information so that I can complete the booking process." "As a user, I want to securely checkout and complete my purchase so that I can receive my order."
Set-ParagraphText "As a traveler, I want to receive a confirmation email with my booking details after completing the purchase so that I have a record of my flight reservation." "As an administrator, I want to manage product information and inventory so that I can ensure accurate product listings."

# 2) Insert a new bullet paragraph right after the just-updated "administrator" bullet.
# InsertParagraphAfter() on a ListBullet paragraph naturally inherits the
# ListBullet style, so there's no need to (re)apply a style explicitly.
$target = Get-ParagraphByText "As an administrator, I want to manage product information and inventory so that I can ensure accurate product listings."
$target.Range.InsertParagraphAfter() | Out-Null
$newPara = $target.Next()
$newPara.Range.Text = "As an administrator, I want to track and manage orders so that I can fulfill customer requests and monitor sales."

# 3) Remove the "Technical Requirements" bullet list (Backend/Frontend/Database/Deployment)
Remove-ParagraphByText "Backend: Python/Django"
Remove-ParagraphByText "Frontend: React.js"
Remove-ParagraphByText "Database: PostgreSQL"
Remove-ParagraphByText "Deployment: Docker/Kubernetes on cloud platform"

# 4) Remove the "Assumptions" bullet list (Payment service / Admin user role)
Remove-ParagraphByText "Payment service is integrated"
Remove-ParagraphByText "Admin user role has access to all booking data"

# 5) Remove the "Open Questions / Risks" bullet list (Query performance / Integration with flight seat map)
Remove-ParagraphByText "Query performance under heavy load"
Remove-ParagraphByText "Integration with flight seat map"
